# Chapter 3.2.4 "Calculating median and variance" -- add the intro
# paragraph runs, flip two empty paragraphs' rFonts hint to eastAsia,
# insert one more blank spacer paragraph, and mark the page break that
# now lands on the "3.3 TODO" heading run.

$d = $word.ActiveDocument

# Locate the "3.2.4 ..." heading paragraph, then walk forward to the
# three short paragraphs that follow it (the bookmark paragraph, the
# spacer, and the second spacer) and finally the "3.3 TODO" heading.
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*3.2.4*") {
        $heading = $cand
        break
    }
}
if ($heading -eq $null) {
    throw "could not locate the 3.2.4 heading paragraph"
}

$pBookmark = $heading.Next()
$pSpacer1  = $pBookmark.Next()
$pSpacer2  = $pSpacer1.Next()
$pNextHeading = $pSpacer2.Next()

# 1) Paragraph that carries the "_GoBack" bookmark: give its pPr rFonts
#    an eastAsia hint and insert the seven new runs of body text before
#    the bookmark.
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="52B62399" w14:textId="308092CE" w:rsidR="00342F2F" w:rsidRPr="00133633" w:rsidRDefault="00342F2F" w:rsidP="00FB44F4"><w:pPr><w:ind w:firstLineChars="200" w:firstLine="480"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>本节将介绍如何计算股票价格的中位数和方差。</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>由于中位数和方差计算所需的原始数据——收盘价，已经在上一节中获取并保存在</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>lose</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>中，因此我们可以直接使用</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>close</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>中的数据。</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>中位数的计算有两种方法：</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pBookmark.Range.InsertXML($xml1)

# 2) The following (empty) paragraph: same eastAsia hint flip on rFonts.
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4F8CA3BB" w14:textId="61720DEE" w:rsidR="00342F2F" w:rsidRDefault="00342F2F" w:rsidP="00FB44F4"><w:pPr><w:ind w:firstLineChars="200" w:firstLine="480"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pSpacer1.Range.InsertXML($xml2)

# 3) Insert a brand-new blank paragraph right after the second spacer
#    paragraph (before the "3.3 TODO" heading).
$xml3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:ind w:firstLineChars="200" w:firstLine="480"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint = $d.Range($pNextHeading.Range.Start, $pNextHeading.Range.Start)
$insertionPoint.InsertXML($xml3)

# Re-resolve the "3.3 TODO" heading paragraph (the package above shifted
# ranges) and stamp a lastRenderedPageBreak inside its first run.
$pNextHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*TODO*") {
        $pNextHeading = $cand
        break
    }
}
if ($pNextHeading -eq $null) {
    throw "could not locate the 3.3 TODO heading paragraph"
}
$xml4 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="53422AAD" w14:textId="5DA8D8AC" w:rsidR="00342F2F" w:rsidRPr="008A0341" w:rsidRDefault="00342F2F" w:rsidP="00342F2F"><w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="008A0341"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>3</w:t></w:r><w:r w:rsidRPr="008A0341"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidRPr="008A0341"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00264400"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>TODO:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>绘制分析数据统计图（七，八</w:t></w:r><w:r w:rsidR="00E87D87"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>节</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>）</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pNextHeading.Range.InsertXML($xml4)

Write-Host "edit applied"
